$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from E1 onto the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New header labels
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New data row values (boolean FALSE)
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false
